$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3139.000735912728
$ws.Range("C2").Value = 1934.53925876059
$ws.Range("D2").Value = 7601.394031554866
